$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark (it will be re-added at the new location at the end).
$existingBm = $d.Bookmarks.Item("_GoBack")
$existingBm.Delete()

# 2. Append the new paragraphs of report content after the current last paragraph
#    ("Notes: see page source for file details").

# --- paragraph 1/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 2/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 3/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Could do examples of games and accuray when little/some/tons of money is bet  on it to prove my conclusion")

# --- paragraph 4/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 5/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about my conclusion again (the hypothesis for calculating probabiltiies)")

# --- paragraph 6/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 7/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about program design")

# --- paragraph 8/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 9/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about the significance but also insignificance of closing markets on events ")

# --- paragraph 10/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 11/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about viability of markets ")

# --- paragraph 12/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 13/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about what I identified as good markets")

# --- paragraph 14/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 15/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about possible games that are trackable and events")

# --- paragraph 16/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 17/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about issues (modelling)")

# --- paragraph 18/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 19/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Talk about current modelling strategy (x pointsand calculate gradient and ensure that the pattern follows")
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore(" (from graph data)")

# --- paragraph 20/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()

# --- paragraph 21/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Trying to extract as much data for free from betfair")

# --- paragraph 22/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("UI isn’t a big deal but taking inspiration from football manager")

# --- paragraph 23/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Trying to cater for different sports (for football x vs y means x is home")

# --- paragraph 24/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Basketball is x @ y so y is home, need to distinguish")

# --- paragraph 25/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("There’s some free data you can get but the way it’s extracted differs slightly per sport")

# --- paragraph 26/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Some sports are totally unviable for obvious reasons, we work with the hypothesis that more market activity = more reactive to match events although there’s the chance that games are totally stagnant so it means little at certain times.")

# --- paragraph 27/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Evaluation")

# --- paragraph 28/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Test code written in the project so that all json replies that are input are saved, this allows reserialization and refeeding the program very quickly for expected output")

# --- paragraph 29/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Deciding possible events and markets because tracking of all market matched/unmatched")

# --- paragraph 30/30 ---
$insPt = $d.Range($d.Content.End, $d.Content.End)
$insPt.InsertParagraphAfter()
$t = $d.Range($d.Content.End, $d.Content.End)
$t.InsertBefore("Chapter onbetfair terminology and api structure")

# Place the _GoBack bookmark collapsed right after the text of this final paragraph
# (matching the original pattern of a zero-width bookmark right before the paragraph
# mark). Bookmarks.Add on a truly collapsed range landing exactly on the paragraph-mark
# slot is mishandled by this host, so we insert a one-character placeholder, wrap the
# bookmark around it, then delete the placeholder -- leaving a correctly collapsed
# bookmark behind, the same way Word itself ends up with one after such edits.
$lastPara = $d.Paragraphs.Last
$markerPos = $lastPara.Range.End - 1
$markerRange = $d.Range($markerPos, $markerPos)
$markerRange.InsertBefore("X")
$bmRange = $d.Range($markerPos, $markerPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$cleanupRange = $d.Range($markerPos, $markerPos + 1)
$cleanupRange.Text = ""
